$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.608.63'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.167.28'
$ws.Range('E3').Value = '  -1.87%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.91'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.607'
$ws.Range('E6').Value = '  -2.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '71.44'
$ws.Range('E7').Value = '  -1.04%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.576'
$ws.Range('E9').Value = '  -3.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.00'
$ws.Range('E10').Value = '  -3.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0907'
$ws.Range('E11').Value = '  -3.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.31'
$ws.Range('E12').Value = '  -4.00%  '
$ws.Range('E13').Value = '  -3.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.68'
$ws.Range('E14').Value = '  -3.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.489.80'
$ws.Range('E15').Value = '  -2.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.36'
$ws.Range('E16').Value = '  +1.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.149.21'
$ws.Range('E17').Value = '  -1.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.787'
$ws.Range('E18').Value = '  -5.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.369.90'
$ws.Range('E19').Value = '  -0.95%  '
$ws.Range('E20').Value = '  -3.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '69.59'
$ws.Range('E21').Value = '  -3.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.77'
$ws.Range('E22').Value = '  -5.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.87'
$ws.Range('E23').Value = '  -12.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '225.83'
$ws.Range('E24').Value = '  -1.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.98'
$ws.Range('E25').Value = '  -2.09%  '
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.69'
$ws.Range('E27').Value = '  -5.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.33'
$ws.Range('E28').Value = '  -8.18%  '
$ws.Range('E29').Value = '  -3.30%  '
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '171.18'
$ws.Range('E31').Value = '  +2.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.79'
$ws.Range('E32').Value = '  -2.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.79'
$ws.Range('E33').Value = '  +6.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0769'
$ws.Range('E34').Value = '  -2.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.12'
$ws.Range('E35').Value = '  -8.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.121'
$ws.Range('E36').Value = '  -2.92%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.105'
$ws.Range('E37').Value = '  +0.86%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.30'
$ws.Range('E38').Value = '  +2.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0297'
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.34'
$ws.Range('E40').Value = '  -7.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.06'
$ws.Range('E41').Value = '  -1.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.33'
$ws.Range('E42').Value = '  -4.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '58.63'
$ws.Range('E43').Value = '  -7.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.189'
$ws.Range('E44').Value = '  -3.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.29'
$ws.Range('E45').Value = '  -3.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0963'
$ws.Range('E46').Value = '  -3.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '98.09'
$ws.Range('E47').Value = '  -4.20%  '
$ws.Range('E48').Value = '  -1.58%  '
$ws.Range('E49').Value = '  -3.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.17'
$ws.Range('E50').Value = '  -6.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.62'
$ws.Range('E51').Value = '  -2.61%  '
